# Automatic update of files.
# Rewrites the per-observation data (columns A-J, P, Q, R, Z, AB) for rows
# 2-14 of the "Artfynd" sheet to match the refreshed export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextNumber {
    # Writes a numeric-looking value while keeping the cell's stored type
    # as text (matches the source file, where columns like "Antal" are
    # inline strings such as "10", "60", etc.).
    param($cellRef, [string]$text)
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $text
    $ws.Range($cellRef).NumberFormat = "General"
}

# Row 2
$ws.Range("A2").Value = 111052885
$ws.Range("B2").Value = 89405
$ws.Range("D2").Value = 'NT'
$ws.Range("E2").Value = 1202
$ws.Range("F2").Value = 'Ullticka'
$ws.Range("G2").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H2").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I2").Value = ""
$ws.Range("Q2").Value = 505474.9755723713
$ws.Range("R2").Value = 6913231.449676614
$ws.Range("Z2").Value = '11:52'
$ws.Range("AB2").Value = '11:52'

# Row 3
$ws.Range("A3").Value = 111053811
$ws.Range("B3").Value = 96348
$ws.Range("D3").Value = 'VU'
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = 'Knärot'
$ws.Range("G3").Value = 'Goodyera repens'
$ws.Range("H3").Value = '(L.) R. Br.'
Set-TextNumber "I3" '1'
$ws.Range("J3").Value = 'm²'
$ws.Range("P3").Value = 'Ön Haverö (Ön Haverö), Mpd'
$ws.Range("Q3").Value = 505259.4328842252
$ws.Range("R3").Value = 6913359.747684986

# Row 4
$ws.Range("A4").Value = 111053678
$ws.Range("B4").Value = 96348
$ws.Range("D4").Value = 'VU'
$ws.Range("E4").Value = 220787
$ws.Range("F4").Value = 'Knärot'
$ws.Range("G4").Value = 'Goodyera repens'
$ws.Range("H4").Value = '(L.) R. Br.'
Set-TextNumber "I4" '60'
$ws.Range("J4").Value = 'plantor/tuvor'
$ws.Range("Q4").Value = 505315.2145614849
$ws.Range("R4").Value = 6913377.483213602
$ws.Range("Z4").Value = '12:34'
$ws.Range("AB4").Value = '12:34'

# Row 5
$ws.Range("A5").Value = 111053875
$ws.Range("B5").Value = 96348
$ws.Range("D5").Value = 'VU'
$ws.Range("E5").Value = 220787
$ws.Range("F5").Value = 'Knärot'
$ws.Range("G5").Value = 'Goodyera repens'
$ws.Range("H5").Value = '(L.) R. Br.'
$ws.Range("P5").Value = 'Rörviken (Rörviken), Mpd'
$ws.Range("Q5").Value = 505202.5913005515
$ws.Range("R5").Value = 6913425.139767874

# Row 6
$ws.Range("A6").Value = 111053489
$ws.Range("B6").Value = 96348
Set-TextNumber "I6" '50'
$ws.Range("Q6").Value = 505380.7652265744
$ws.Range("R6").Value = 6913393.377649955

# Row 7
$ws.Range("A7").Value = 111053763
$ws.Range("B7").Value = 96348
$ws.Range("D7").Value = 'VU'
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = 'Knärot'
$ws.Range("G7").Value = 'Goodyera repens'
$ws.Range("H7").Value = '(L.) R. Br.'
Set-TextNumber "I7" '10'
$ws.Range("Q7").Value = 505308.2479573332
$ws.Range("R7").Value = 6913371.434886473
$ws.Range("Z7").Value = '12:34'
$ws.Range("AB7").Value = '12:34'

# Row 8 (only Taxonsorteringsordning changes)
$ws.Range("B8").Value = 96348

# Row 9 (only Taxonsorteringsordning changes)
$ws.Range("B9").Value = 96348

# Row 10
$ws.Range("A10").Value = 111053919
$ws.Range("B10").Value = 89405
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 1202
$ws.Range("F10").Value = 'Ullticka'
$ws.Range("G10").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H10").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I10").Value = ""
$ws.Range("J10").Value = ""
$ws.Range("P10").Value = 'Rörviken (Rörviken), Mpd'
$ws.Range("Q10").Value = 505165.3769719804
$ws.Range("R10").Value = 6913430.654938316

# Row 11
$ws.Range("A11").Value = 111052930
$ws.Range("B11").Value = 96348
$ws.Range("D11").Value = 'VU'
$ws.Range("E11").Value = 220787
$ws.Range("F11").Value = 'Knärot'
$ws.Range("G11").Value = 'Goodyera repens'
$ws.Range("H11").Value = '(L.) R. Br.'
Set-TextNumber "I11" '20'
$ws.Range("J11").Value = 'plantor/tuvor'
$ws.Range("Q11").Value = 505479.575643972
$ws.Range("R11").Value = 6913262.573581941

# Row 12
$ws.Range("A12").Value = 111053802
$ws.Range("B12").Value = 77268
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 228912
$ws.Range("F12").Value = 'Mörk kolflarnlav'
$ws.Range("G12").Value = 'Carbonicola myrmecina'
$ws.Range("H12").Value = '(Ach.) Bendiksby & Timdal'
$ws.Range("P12").Value = 'Ön Haverö (Ön Haverö), Mpd'
$ws.Range("Q12").Value = 505273.4006640643
$ws.Range("R12").Value = 6913350.017071255

# Row 13
$ws.Range("A13").Value = 111053254
$ws.Range("B13").Value = 77268
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = 228912
$ws.Range("F13").Value = 'Mörk kolflarnlav'
$ws.Range("G13").Value = 'Carbonicola myrmecina'
$ws.Range("H13").Value = '(Ach.) Bendiksby & Timdal'
$ws.Range("I13").Value = ""
$ws.Range("J13").Value = ""
$ws.Range("Q13").Value = 505418.1403733313
$ws.Range("R13").Value = 6913288.945539451

# Row 14
$ws.Range("A14").Value = 111053357
$ws.Range("B14").Value = 78604
$ws.Range("D14").Value = 'LC'
$ws.Range("E14").Value = 6461
$ws.Range("F14").Value = 'Norrlandslav'
$ws.Range("G14").Value = 'Nephroma arcticum'
$ws.Range("H14").Value = '(L.) Torss.'
$ws.Range("I14").Value = ""
$ws.Range("J14").Value = ""
$ws.Range("Q14").Value = 505390.5931249987
$ws.Range("R14").Value = 6913355.312167899
$ws.Range("Z14").Value = '11:52'
$ws.Range("AB14").Value = '11:52'
